$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DAMSLTag (column I) and DialogAct (column J) values
# following a re-run of SGNN dialog act annotation.
$updates = @(
    @{ Row = 5; I = 'ba'; J = 'Appreciation' }
    @{ Row = 7; I = 'ba'; J = 'Appreciation' }
    @{ Row = 58; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 74; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 79; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 97; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 138; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 140; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 152; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 166; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 175; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 177; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 186; I = 'aa'; J = 'Agree/Accept' }
    @{ Row = 221; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 235; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 252; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 266; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 273; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 281; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 288; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 300; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 320; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 322; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 330; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 334; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 354; I = 'ba'; J = 'Appreciation' }
    @{ Row = 369; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 377; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 387; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 390; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 396; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 401; I = 'ba'; J = 'Appreciation' }
    @{ Row = 450; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 453; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 457; I = '%'; J = 'Uninterpretable' }
    @{ Row = 464; I = 'ba'; J = 'Appreciation' }
    @{ Row = 469; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 474; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 475; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 484; I = 'b'; J = 'Acknowledge (Backchannel)' }
    @{ Row = 495; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 496; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 500; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 501; I = 'sv'; J = 'Statement-opinion' }
    @{ Row = 519; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 520; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 521; I = 'fc'; J = 'Conventional-closing' }
    @{ Row = 523; I = '%'; J = 'Uninterpretable' }
    @{ Row = 525; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 528; I = 'sd'; J = 'Statement-non-opinion' }
    @{ Row = 530; I = 'sv'; J = 'Statement-opinion' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.I
    $ws.Cells.Item($u.Row, 10).Value = $u.J
}
